$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.087.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "'1.830.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("E4").Value = "  +0.91%  "
$ws.Range("D5").Value = "'312.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").Value = "'1.006"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.81%  "
$ws.Range("D7").Value = "'0.4706"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "'0.3687"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.82%  "
$ws.Range("D9").Value = "'0.07376"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("D10").Value = "'0.8794"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.25%  "
$ws.Range("D11").Value = "'20.36"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("D12").Value = "'1.825.37"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.67%  "
$ws.Range("D13").Value = "'0.07300"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.08%  "
$ws.Range("D14").Value = "'5.456"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.47%  "
$ws.Range("D15").Value = "'92.52"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("D16").Value = "'6.545"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("D18").Value = "'0.000008771"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("D20").Value = "'14.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("D21").Value = "'27.091.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "'5.301"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.14%  "
$ws.Range("D23").Value = "'10.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.68%  "
$ws.Range("D24").Value = "'2.060.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.95%  "
$ws.Range("D25").Value = "'1.896"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").Value = "'151.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").Value = "'18.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").Value = "'2.155"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.11%  "
$ws.Range("D29").Value = "'5.262"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.19%  "
$ws.Range("D30").Value = "'116.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.09%  "
$ws.Range("D31").Value = "'0.08914"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.52%  "
$ws.Range("D32").Value = "'0.7579"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.04%  "
$ws.Range("D33").Value = "'1.167"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").Value = "'4.527"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.79%  "
$ws.Range("D35").Value = "'2.926"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("E36").Value = "  +0.84%  "
$ws.Range("D37").Value = "'1.099"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").Value = "'0.05329"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("D39").Value = "'0.01954"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").Value = "'2.996"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.30%  "
$ws.Range("D41").Value = "'2.416"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("D42").Value = "'7.272"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").Value = "'0.5327"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.12%  "
$ws.Range("D44").Value = "'0.1662"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.42%  "
$ws.Range("D45").Value = "'8.547"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").Value = "'0.4944"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.89%  "
$ws.Range("D47").Value = "'10.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.65%  "
$ws.Range("D48").Value = "'1.006"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.83%  "
$ws.Range("D49").Value = "'1.667"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.16%  "
$ws.Range("D50").Value = "'103.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("D51").Value = "'0.06309"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.10%  "
